$d = $word.ActiveDocument

function Get-ParaByText($doc, $targetText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $targetText) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Change 1: "UC04-App Read Message" -> "UC04-App Læs Beskeder"
#   ("Read Message" becomes "Beskeder"; " Læs" is added before the
#    existing space that separates "-App" from the label.)
# ---------------------------------------------------------------------
$p1 = Get-ParaByText $d "UC04-App Read Message"
if ($p1 -ne $null) {
    $pr1 = $p1.Range
    $sr1 = $d.Range($pr1.Start, $pr1.End)
    $sr1.Find.Execute("Read Message", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
    $rmRange = $d.Range($sr1.Start, $sr1.End)
    $rmRange.Text = "Beskeder"

    $p1b = Get-ParaByText $d "UC04-App Beskeder"
    if ($p1b -ne $null) {
        $pr1b = $p1b.Range
        $sr1b = $d.Range($pr1b.Start, $pr1b.End)
        $sr1b.Find.Execute(" Beskeder", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
        $insPoint1 = $d.Range($sr1b.Start, $sr1b.Start)
        $insPoint1.InsertBefore(" Læs")
    }
}

# ---------------------------------------------------------------------
# Change 2: "Fysioapp viser klientens indbakke." ->
#           "Frederiksberg Sportsklinik viser klientens indbakke."
#   (Both resulting runs keep the existing Arial formatting.)
# ---------------------------------------------------------------------
$p2 = Get-ParaByText $d "Fysioapp viser klientens indbakke."
if ($p2 -ne $null) {
    $pr2 = $p2.Range
    $sr2 = $d.Range($pr2.Start, $pr2.End)
    $sr2.Find.Execute("Fysioapp", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
    $rng2 = $d.Range($sr2.Start, $sr2.End)
    $rng2.Text = "Frederiksberg Sportsklinik"
}

# ---------------------------------------------------------------------
# Change 3: "Fysioapp viser beskeden for klienten." ->
#           "Frederiksberg Sportsklinik viser beskeden for klienten."
#   (Here "Frederiksberg Sportsklinik" becomes its own run carrying no
#    direct character formatting, while the rest of the sentence keeps
#    the original da-DK run formatting.)
# ---------------------------------------------------------------------
$p3 = Get-ParaByText $d "Fysioapp viser beskeden for klienten."
if ($p3 -ne $null) {
    # Insert the new name as an unformatted run at the very start of the
    # paragraph FIRST (collapsing the paragraph's own Range avoids
    # inheriting the following run's formatting). Doing this before the
    # "Fysioapp" removal keeps the new run from being reabsorbed into
    # the surrounding lang=da-DK formatting.
    $ins3 = $p3.Range
    $ins3.Collapse(1)
    $ins3.InsertBefore("Frederiksberg Sportsklinik")

    # Now remove the old "Fysioapp" word, leaving the trailing space
    # behind so the sentence reads correctly.
    $p3b = Get-ParaByText $d "Frederiksberg SportsklinikFysioapp viser beskeden for klienten."
    if ($p3b -ne $null) {
        $pr3b = $p3b.Range
        $sr3b = $d.Range($pr3b.Start, $pr3b.End)
        $sr3b.Find.Execute("Fysioapp", $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
        $rng3b = $d.Range($sr3b.Start, $sr3b.End)
        $rng3b.Text = ""
    }
}
